# --- Fixed update to Excel issue ---
# 1) Rename the "Requested quantity" header on the two existing sheets.
# 2) Add a new "PO Forecast" sheet (built from a copy of "Weekly Quantity" so
#    it inherits identical header / date-column formatting) with the
#    forecasted PO quantities + confidence interval.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Range("B1").Value = "Weekly_PO_Qty"

$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# Build the new sheet by copying an existing one (after the last sheet) so it
# starts out with the same header/date styles already present in the style
# table instead of minting brand-new font/style entries.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws1.Copy($null, $lastSheet)
$ws3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3.Name = "PO Forecast"

# Headers
$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

# C1/D1 are brand-new cells on this sheet, so pick up the same bold/bordered
# header formatting already used by A1/B1 via a format-only paste.
$ws3.Range("A1").Copy()
$ws3.Range("C1:D1").PasteSpecial(-4122)

# The copied sheet only has 22 data rows (rows 2-23); extend the date-style
# formatting on column A down to row 30 before writing values.
$ws3.Range("A22").Copy()
$ws3.Range("A23:A30").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$data = @(
    @(44934.99999999999, 127, -8.902605420146568, 256.8313745878708),
    @(44955.99999999999, 120, -22.07776844778072, 251.6170328771354),
    @(44969.99999999999, 116, -16.72188997173397, 251.2273295891223),
    @(44990.99999999999, 109, -24.64982909330645, 237.9964052993855),
    @(44997.99999999999, 107, -30.18185836465519, 235.6480024132763),
    @(45011.99999999999, 102, -29.1335198083676, 244.4308499667835),
    @(45018.99999999999, 100, -35.36726837243619, 236.6250875138792),
    @(45025.99999999999, 98, -25.2680619006217, 235.7070225684175),
    @(45032.99999999999, 95, -38.81645266744457, 227.3263439135131),
    @(45039.99999999999, 93, -46.6707384009835, 223.9582649387747),
    @(45046.99999999999, 91, -34.22111542344957, 217.3151777900927),
    @(45053.99999999999, 89, -39.5060730938435, 215.7219022334855),
    @(45074.99999999999, 82, -54.27207212895216, 209.5393923527178),
    @(45081.99999999999, 80, -53.84174322506442, 210.9741593772221),
    @(45095.99999999999, 75, -56.75567324508828, 205.2473565139905),
    @(45102.99999999999, 73, -62.05077984564738, 204.5550789737555),
    @(45109.99999999999, 71, -52.56495535699972, 198.6883122194039),
    @(45116.99999999999, 68, -69.70979955949528, 192.3961033176145),
    @(45123.99999999999, 66, -65.99393550383954, 204.8613394432608),
    @(45130.99999999999, 64, -72.3199129503754, 195.5783892310787),
    @(45137.99999999999, 62, -70.44422296973053, 191.2928743628627),
    @(45144.99999999999, 59, -74.57857450277049, 191.9589042131466),
    @(45151.99999999999, 57, -73.72119256965323, 188.986607242907),
    @(45158.99999999999, 55, -73.7412261785429, 189.9169293600306),
    @(45165.99999999999, 53, -81.66152946979408, 180.8539388696355),
    @(45172.99999999999, 50, -84.43410448755832, 167.2820714405042),
    @(45179.99999999999, 48, -87.98817554423913, 185.6988720227306),
    @(45186.99999999999, 46, -92.76796452735709, 180.4970192618647),
    @(45193.99999999999, 44, -96.02216368322134, 171.4504372642787)
)

$r = 2
foreach ($row in $data) {
    $ws3.Cells.Item($r, 1).Value = $row[0]
    $ws3.Cells.Item($r, 2).Value = $row[1]
    $ws3.Cells.Item($r, 3).Value = $row[2]
    $ws3.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

[void]$ws3.Range("A1").Select()
